$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 234 ("「知っていますか？」..." quiz entry) entirely.
# This shifts all subsequent rows (235..393) up by one, matching the
# target dimension shrinking from A1:C393 to A1:C392.
$ws.Rows.Item(234).Delete()
